$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 554.7143
$ws.Range("I41").Value = 344.75
$ws.Range("J41").Value = 638.7
$ws.Range("K41").Value = 344.75
$ws.Range("L41").Value = 638.7
$ws.Range("M41").Value = 95.25
$ws.Range("N41").Value = -1518.7
# Row 64
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 5000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -5496
# Row 67
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = 5000
$ws.Range("N67").Value = -6716
# Row 132
$ws.Range("H132").Value = 479924.38
$ws.Range("I132").Value = 289343.38
$ws.Range("J132").Value = 1432829.4
$ws.Range("K132").Value = 868030.14
$ws.Range("L132").Value = 4298488.199999999
$ws.Range("M132").Value = -865500.14
$ws.Range("N132").Value = -4303548.199999999
# Row 138
$ws.Range("H138").Value = 1709.34
$ws.Range("I138").Value = 598.6923
$ws.Range("J138").Value = 2419.4263
$ws.Range("K138").Value = 1796.0769
$ws.Range("L138").Value = 7258.2789
$ws.Range("M138").Value = 3343.9231
$ws.Range("N138").Value = -17538.2789

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5003.3384
$ws.Range("I32").Value = 4429.0376
$ws.Range("J32").Value = 7032.533
$ws.Range("K32").Value = 4429.0376
$ws.Range("L32").Value = 7032.533
$ws.Range("M32").Value = -4142.0376
$ws.Range("N32").Value = -7606.533

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2442.4
$ws.Range("I86").Value = 2083.4
$ws.Range("J86").Value = 2801.4
$ws.Range("K86").Value = 2083.4
$ws.Range("L86").Value = 2801.4
$ws.Range("M86").Value = -960.4000000000001
$ws.Range("N86").Value = -5047.4
# Row 89
$ws.Range("H89").Value = 2442.4
$ws.Range("I89").Value = 2083.4
$ws.Range("J89").Value = 2801.4
$ws.Range("K89").Value = 10417
$ws.Range("L89").Value = 14007
$ws.Range("M89").Value = -4801
$ws.Range("N89").Value = -25239
# Row 122
$ws.Range("H122").Value = 41781.43
$ws.Range("J122").Value = 41781.43
$ws.Range("L122").Value = 41781.43
$ws.Range("N122").Value = -51581.43
# Row 134
$ws.Range("H134").Value = 1534.0625
$ws.Range("I134").Value = 940.7059
$ws.Range("J134").Value = 2975.0715
$ws.Range("K134").Value = 2822.1177
$ws.Range("L134").Value = 8925.2145
$ws.Range("M134").Value = -287.1177000000002
$ws.Range("N134").Value = -13995.2145

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1779.6794
$ws.Range("I58").Value = 1588.0625
$ws.Range("J58").Value = 2655.6428
$ws.Range("K58").Value = 1588.0625
$ws.Range("L58").Value = 2655.6428
$ws.Range("M58").Value = -1385.0625
$ws.Range("N58").Value = -3061.6428
# Row 134
$ws.Range("H134").Value = 3005.423
$ws.Range("I134").Value = 3294.3242
$ws.Range("J134").Value = 2292.8
$ws.Range("K134").Value = 9882.972600000001
$ws.Range("L134").Value = 6878.400000000001
$ws.Range("M134").Value = -7347.972600000001
$ws.Range("N134").Value = -11948.4
# Row 136
$ws.Range("H136").Value = 1779.6794
$ws.Range("I136").Value = 1588.0625
$ws.Range("J136").Value = 2655.6428
$ws.Range("K136").Value = 4764.1875
$ws.Range("L136").Value = 7966.928400000001
$ws.Range("M136").Value = -2214.1875
$ws.Range("N136").Value = -13066.9284

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 788.1852
$ws.Range("I113").Value = 664.5
$ws.Range("J113").Value = 887.13336
$ws.Range("K113").Value = 1993.5
$ws.Range("L113").Value = 2661.40008
$ws.Range("M113").Value = 176.5
$ws.Range("N113").Value = -7001.40008
# Row 131
$ws.Range("H131").Value = 6667494.5
$ws.Range("I131").Value = 100000280
$ws.Range("J131").Value = 867.02856
$ws.Range("K131").Value = 300000840
$ws.Range("L131").Value = 2601.08568
$ws.Range("M131").Value = -299995800
$ws.Range("N131").Value = -12681.08568
# Row 137
$ws.Range("H137").Value = 3557.3635
$ws.Range("J137").Value = 4328.875
$ws.Range("L137").Value = 12986.625
$ws.Range("N137").Value = -23186.625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 83335300
$ws.Range("I80").Value = 125001450
$ws.Range("K80").Value = 125001450
$ws.Range("M80").Value = -125000452
# Row 83
$ws.Range("H83").Value = 83335300
$ws.Range("I83").Value = 125001450
$ws.Range("K83").Value = 625007250
$ws.Range("M83").Value = -625002258
# Row 126
$ws.Range("H126").Value = 3352.32
$ws.Range("I126").Value = 3004.2026
$ws.Range("J126").Value = 4661.905
$ws.Range("K126").Value = 9012.6078
$ws.Range("L126").Value = 13985.715
$ws.Range("M126").Value = -6542.6078
$ws.Range("N126").Value = -18925.715

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 26
$ws.Range("H26").Value = 29991
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 29991
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 29991
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -30581
# Row 45
$ws.Range("H45").Value = 29992.75
$ws.Range("I45").Value = 20040.5
$ws.Range("J45").Value = 39945
$ws.Range("K45").Value = 20040.5
$ws.Range("L45").Value = 39945
$ws.Range("M45").Value = -19633.5
$ws.Range("N45").Value = -40759
# Row 74
$ws.Range("H74").Value = 45843.4
$ws.Range("J74").Value = 45843.4
$ws.Range("L74").Value = 45843.4
$ws.Range("N74").Value = -47839.4
# Row 77
$ws.Range("H77").Value = 45843.4
$ws.Range("J77").Value = 45843.4
$ws.Range("L77").Value = 137530.2
$ws.Range("N77").Value = -147514.2
# Row 82
$ws.Range("H82").Value = 1478.56
$ws.Range("I82").Value = 651.2727
$ws.Range("J82").Value = 2128.5715
$ws.Range("K82").Value = 651.2727
$ws.Range("L82").Value = 2128.5715
$ws.Range("M82").Value = -290.2727
$ws.Range("N82").Value = -2850.5715
# Row 85
$ws.Range("H85").Value = 1478.56
$ws.Range("I85").Value = 651.2727
$ws.Range("J85").Value = 2128.5715
$ws.Range("K85").Value = 651.2727
$ws.Range("L85").Value = 2128.5715
$ws.Range("M85").Value = 596.7273
$ws.Range("N85").Value = -4624.5715
# Row 107
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
# Row 122
$ws.Range("H122").Value = 3914.4211
$ws.Range("I122").Value = 2019.5
$ws.Range("J122").Value = 6019.8887
$ws.Range("K122").Value = 6058.5
$ws.Range("L122").Value = 18059.6661
$ws.Range("M122").Value = -3608.5
$ws.Range("N122").Value = -22959.6661
# Row 136
$ws.Range("H136").Value = 2556.12
$ws.Range("I136").Value = 929.58826
$ws.Range("K136").Value = 2788.76478
$ws.Range("M136").Value = -238.76478

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5377514.5
$ws.Range("I132").Value = 712.9231
$ws.Range("J132").Value = 33336884
$ws.Range("K132").Value = 2138.7693
$ws.Range("L132").Value = 100010652
$ws.Range("M132").Value = 391.2307000000001
$ws.Range("N132").Value = -100015712

